$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Copy Sheet1's data (A1:F21) into Sheet2 starting at row 13, but with the
# columns rotated: Sheet1 column A (Station Name) becomes Sheet2 column F,
# and Sheet1 columns B:F shift left into Sheet2 columns A:E.
$ws1.Range("B1:F21").Copy()
$ws2.Range("A13").PasteSpecial()

$ws1.Range("A1:A21").Copy()
$ws2.Range("F13").PasteSpecial()

# Re-apply the same cell formatting (style) that Sheet1 used, onto the newly
# populated Sheet2 range.
$ws1.Range("A1:F21").Copy()
$ws2.Range("A13:F33").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Update the selections left on each sheet.
$ws1.Range("A1:F21").Select()
$ws2.Range("I20").Select()

# Sheet2 is now the active/visible tab.
$ws2.Activate()
